$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2805039  # H98
$ws.Cells.Item(98, 9).Value = 31053.4  # I98
$ws.Cells.Item(98, 10).Value = 22222938  # J98
$ws.Cells.Item(98, 11).Value = 31053.4  # K98
$ws.Cells.Item(98, 12).Value = 22222938  # L98
$ws.Cells.Item(98, 13).Value = -29555.4  # M98
$ws.Cells.Item(98, 14).Value = -22225934  # N98
$ws.Cells.Item(122, 8).Value = 2805039  # H122
$ws.Cells.Item(122, 9).Value = 31053.4  # I122
$ws.Cells.Item(122, 10).Value = 22222938  # J122
$ws.Cells.Item(122, 11).Value = 93160.20000000001  # K122
$ws.Cells.Item(122, 12).Value = 66668814  # L122
$ws.Cells.Item(122, 13).Value = -90710.20000000001  # M122
$ws.Cells.Item(122, 14).Value = -66673714  # N122
$ws.Cells.Item(140, 8).Value = 61646.875  # H140
$ws.Cells.Item(140, 9).Value = 10000  # I140
$ws.Cells.Item(140, 11).Value = 10000  # K140
$ws.Cells.Item(140, 13).Value = -4820  # M140

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(60, 8).Value = 1000  # H60
$ws.Cells.Item(60, 9).Value = 1000  # I60
$ws.Cells.Item(60, 11).Value = 1000  # K60
$ws.Cells.Item(60, 13).Value = -267  # M60
$ws.Cells.Item(61, 8).Value = 7453.6  # H61
$ws.Cells.Item(61, 9).Value = 7835.778  # I61
$ws.Cells.Item(61, 10).Value = 4014  # J61
$ws.Cells.Item(61, 11).Value = 7835.778  # K61
$ws.Cells.Item(61, 12).Value = 4014  # L61
$ws.Cells.Item(61, 13).Value = -7623.778  # M61
$ws.Cells.Item(61, 14).Value = -4438  # N61
$ws.Cells.Item(74, 8).Value = 1291.1666  # H74
$ws.Cells.Item(74, 9).Value = 839.0769  # I74
$ws.Cells.Item(74, 10).Value = 1825.4546  # J74
$ws.Cells.Item(74, 11).Value = 839.0769  # K74
$ws.Cells.Item(74, 12).Value = 1825.4546  # L74
$ws.Cells.Item(74, 13).Value = 34.92309999999998  # M74
$ws.Cells.Item(74, 14).Value = -3573.4546  # N74
$ws.Cells.Item(77, 8).Value = 1291.1666  # H77
$ws.Cells.Item(77, 9).Value = 839.0769  # I77
$ws.Cells.Item(77, 10).Value = 1825.4546  # J77
$ws.Cells.Item(77, 11).Value = 4195.3845  # K77
$ws.Cells.Item(77, 12).Value = 9127.273000000001  # L77
$ws.Cells.Item(77, 13).Value = 172.6154999999999  # M77
$ws.Cells.Item(77, 14).Value = -17863.273  # N77
$ws.Cells.Item(88, 8).Value = 5903  # H88
$ws.Cells.Item(88, 10).Value = 5903  # J88
$ws.Cells.Item(88, 12).Value = 5903  # L88
$ws.Cells.Item(88, 14).Value = -6715  # N88
$ws.Cells.Item(91, 8).Value = 5903  # H91
$ws.Cells.Item(91, 10).Value = 5903  # J91
$ws.Cells.Item(91, 12).Value = 5903  # L91
$ws.Cells.Item(91, 14).Value = -8711  # N91
$ws.Cells.Item(136, 8).Value = 7453.6  # H136
$ws.Cells.Item(136, 9).Value = 7835.778  # I136
$ws.Cells.Item(136, 10).Value = 4014  # J136
$ws.Cells.Item(136, 11).Value = 23507.334  # K136
$ws.Cells.Item(136, 12).Value = 12042  # L136
$ws.Cells.Item(136, 13).Value = -20957.334  # M136
$ws.Cells.Item(136, 14).Value = -17142  # N136
$ws.Cells.Item(138, 8).Value = 63933.332  # H138
$ws.Cells.Item(138, 10).Value = 63933.332  # J138
$ws.Cells.Item(138, 12).Value = 63933.332  # L138
$ws.Cells.Item(138, 14).Value = -74213.33199999999  # N138
$ws.Cells.Item(139, 8).Value = 65000  # H139
$ws.Cells.Item(139, 10).Value = 65000  # J139
$ws.Cells.Item(139, 12).Value = 65000  # L139
$ws.Cells.Item(139, 14).Value = -75280  # N139

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 7989.25  # H75
$ws.Cells.Item(75, 9).Value = 3982.8  # I75
$ws.Cells.Item(75, 10).Value = 14666.667  # J75
$ws.Cells.Item(75, 11).Value = 3982.8  # K75
$ws.Cells.Item(75, 12).Value = 14666.667  # L75
$ws.Cells.Item(75, 13).Value = -3046.8  # M75
$ws.Cells.Item(75, 14).Value = -16538.667  # N75
$ws.Cells.Item(78, 8).Value = 7989.25  # H78
$ws.Cells.Item(78, 9).Value = 3982.8  # I78
$ws.Cells.Item(78, 10).Value = 14666.667  # J78
$ws.Cells.Item(78, 11).Value = 11948.4  # K78
$ws.Cells.Item(78, 12).Value = 44000.001  # L78
$ws.Cells.Item(78, 13).Value = -7268.400000000001  # M78
$ws.Cells.Item(78, 14).Value = -53360.001  # N78
$ws.Cells.Item(132, 8).Value = 35000  # H132
$ws.Cells.Item(132, 10).Value = 35000  # J132
$ws.Cells.Item(132, 12).Value = 35000  # L132
$ws.Cells.Item(132, 14).Value = -45120  # N132
$ws.Cells.Item(138, 8).Value = 63266.668  # H138
$ws.Cells.Item(138, 10).Value = 63266.668  # J138
$ws.Cells.Item(138, 12).Value = 63266.668  # L138
$ws.Cells.Item(138, 14).Value = -73546.66800000001  # N138
$ws.Cells.Item(140, 8).Value = 89900  # H140
$ws.Cells.Item(140, 10).Value = 89900  # J140
$ws.Cells.Item(140, 12).Value = 89900  # L140
$ws.Cells.Item(140, 14).Value = -100260  # N140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1483  # H58
$ws.Cells.Item(58, 9).Value = 897.1111  # I58
$ws.Cells.Item(58, 10).Value = 2010.3  # J58
$ws.Cells.Item(58, 11).Value = 897.1111  # K58
$ws.Cells.Item(58, 12).Value = 2010.3  # L58
$ws.Cells.Item(58, 13).Value = -694.1111  # M58
$ws.Cells.Item(58, 14).Value = -2416.3  # N58
$ws.Cells.Item(64, 8).Value = 0  # H64
$ws.Cells.Item(64, 10).Value = 0  # J64
$ws.Cells.Item(64, 12).Value = 0  # L64
$ws.Cells.Item(67, 8).Value = 0  # H67
$ws.Cells.Item(67, 10).Value = 0  # J67
$ws.Cells.Item(67, 12).Value = 0  # L67
$ws.Cells.Item(136, 8).Value = 1483  # H136
$ws.Cells.Item(136, 9).Value = 897.1111  # I136
$ws.Cells.Item(136, 10).Value = 2010.3  # J136
$ws.Cells.Item(136, 11).Value = 2691.3333  # K136
$ws.Cells.Item(136, 12).Value = 6030.9  # L136
$ws.Cells.Item(136, 13).Value = -141.3332999999998  # M136
$ws.Cells.Item(136, 14).Value = -11130.9  # N136
$ws.Cells.Item(140, 8).Value = 88369.25  # H140
$ws.Cells.Item(140, 10).Value = 88369.25  # J140
$ws.Cells.Item(140, 12).Value = 88369.25  # L140
$ws.Cells.Item(140, 14).Value = -98729.25  # N140
$ws.Cells.Item(64, 14).ClearContents()  # N64 (removed)
$ws.Cells.Item(67, 14).ClearContents()  # N67 (removed)

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1002.8947  # H5
$ws.Cells.Item(5, 9).Value = 601.8421  # I5
$ws.Cells.Item(5, 10).Value = 1403.9474  # J5
$ws.Cells.Item(5, 11).Value = 1805.5263  # K5
$ws.Cells.Item(5, 12).Value = 4211.8422  # L5
$ws.Cells.Item(5, 13).Value = -1693.5263  # M5
$ws.Cells.Item(5, 14).Value = -4435.8422  # N5
$ws.Cells.Item(135, 8).Value = 1002.8947  # H135
$ws.Cells.Item(135, 9).Value = 601.8421  # I135
$ws.Cells.Item(135, 10).Value = 1403.9474  # J135
$ws.Cells.Item(135, 11).Value = 5416.5789  # K135
$ws.Cells.Item(135, 12).Value = 12635.5266  # L135
$ws.Cells.Item(135, 13).Value = -2881.5789  # M135
$ws.Cells.Item(135, 14).Value = -17705.5266  # N135

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 54290  # H133
$ws.Cells.Item(133, 10).Value = 54290  # J133
$ws.Cells.Item(133, 12).Value = 54290  # L133
$ws.Cells.Item(133, 14).Value = -64410  # N133
$ws.Cells.Item(138, 8).Value = 68207.14  # H138
$ws.Cells.Item(138, 10).Value = 68207.14  # J138
$ws.Cells.Item(138, 12).Value = 68207.14  # L138
$ws.Cells.Item(138, 14).Value = -78487.14  # N138
$ws.Cells.Item(140, 8).Value = 89863  # H140
$ws.Cells.Item(140, 10).Value = 89863  # J140
$ws.Cells.Item(140, 12).Value = 89863  # L140
$ws.Cells.Item(140, 14).Value = -100223  # N140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(62, 8).Value = 0  # H62
$ws.Cells.Item(62, 10).Value = 0  # J62
$ws.Cells.Item(62, 12).Value = 0  # L62
$ws.Cells.Item(65, 8).Value = 0  # H65
$ws.Cells.Item(65, 10).Value = 0  # J65
$ws.Cells.Item(65, 12).Value = 0  # L65
$ws.Cells.Item(133, 8).Value = 86561.734  # H133
$ws.Cells.Item(133, 10).Value = 86561.734  # J133
$ws.Cells.Item(133, 12).Value = 86561.734  # L133
$ws.Cells.Item(133, 14).Value = -91621.734  # N133
$ws.Cells.Item(139, 8).Value = 79475  # H139
$ws.Cells.Item(139, 10).Value = 79475  # J139
$ws.Cells.Item(139, 12).Value = 79475  # L139
$ws.Cells.Item(139, 14).Value = -89755  # N139
$ws.Cells.Item(62, 14).ClearContents()  # N62 (removed)
$ws.Cells.Item(65, 14).ClearContents()  # N65 (removed)

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(82, 8).Value = 30250  # H82
$ws.Cells.Item(82, 10).Value = 30250  # J82
$ws.Cells.Item(82, 12).Value = 30250  # L82
$ws.Cells.Item(82, 14).Value = -31016  # N82
$ws.Cells.Item(85, 8).Value = 30250  # H85
$ws.Cells.Item(85, 10).Value = 30250  # J85
$ws.Cells.Item(85, 12).Value = 30250  # L85
$ws.Cells.Item(85, 14).Value = -32902  # N85
$ws.Cells.Item(133, 8).Value = 19161.25  # H133
$ws.Cells.Item(133, 10).Value = 19161.25  # J133
$ws.Cells.Item(133, 12).Value = 19161.25  # L133
$ws.Cells.Item(133, 14).Value = -29281.25  # N133
$ws.Cells.Item(136, 8).Value = 1113.3572  # H136
$ws.Cells.Item(136, 9).Value = 1134  # I136
$ws.Cells.Item(136, 10).Value = 1051.4286  # J136
$ws.Cells.Item(136, 11).Value = 3402  # K136
$ws.Cells.Item(136, 12).Value = 3154.2858  # L136
$ws.Cells.Item(136, 13).Value = -852  # M136
$ws.Cells.Item(136, 14).Value = -8254.2858  # N136
$ws.Cells.Item(138, 8).Value = 69180  # H138
$ws.Cells.Item(138, 10).Value = 69180  # J138
$ws.Cells.Item(138, 12).Value = 69180  # L138
$ws.Cells.Item(138, 14).Value = -79460  # N138
$ws.Cells.Item(139, 8).Value = 57328.75  # H139
$ws.Cells.Item(139, 10).Value = 57328.75  # J139
$ws.Cells.Item(139, 12).Value = 57328.75  # L139
$ws.Cells.Item(139, 14).Value = -67608.75  # N139
